$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.05284710587315838
$ws.Range("E2").Value = 0.02550090984167503
$ws.Range("G2").Value = -0.04502025361393325
$ws.Range("H2").Value = -0.5352988047422513
$ws.Range("I2").Value = 0.1490515468083845
$ws.Range("J2").Value = 0.1371387993843627
$ws.Range("K2").Value = 0.200626606284561
$ws.Range("L2").Value = 0.1079330116962531
$ws.Range("M2").Value = 0.129293823016267
$ws.Range("N2").Value = -0.5502255047420385
$ws.Range("O2").Value = 0.2537685189959536
$ws.Range("D3").Value = 0.1183818032404943
$ws.Range("E3").Value = -0.0427249620654119
$ws.Range("G3").Value = -0.01940981040636119
$ws.Range("H3").Value = 0.1686289861464561
$ws.Range("I3").Value = 0.1051672521303078
$ws.Range("J3").Value = -0.01637347079980878
$ws.Range("K3").Value = -0.07889299984295564
$ws.Range("L3").Value = 0.1237571271336534
$ws.Range("M3").Value = -0.03019904530512181
$ws.Range("N3").Value = 0.1411896716252118
$ws.Range("O3").Value = 0.02761189664673756
$ws.Range("B4").Value = 0.05284710587315838
$ws.Range("C4").Value = 0.1183818032404943
$ws.Range("E4").Value = 0.2862096796565572
$ws.Range("G4").Value = 0.2858119586419813
$ws.Range("H4").Value = 0.3208449815547014
$ws.Range("I4").Value = -0.03615914419776285
$ws.Range("J4").Value = 0.8042684740305435
$ws.Range("K4").Value = 0.1010925659905756
$ws.Range("L4").Value = 0.1766723565517075
$ws.Range("M4").Value = 0.06202386746339025
$ws.Range("N4").Value = 0.2608281976605054
$ws.Range("O4").Value = -0.09255504229564779
$ws.Range("B5").Value = 0.02550090984167503
$ws.Range("C5").Value = -0.0427249620654119
$ws.Range("D5").Value = 0.2862096796565572
$ws.Range("G5").Value = 0.9628637223458245
$ws.Range("H5").Value = 0.1878429442749909
$ws.Range("I5").Value = 0.4840581834101727
$ws.Range("J5").Value = 0.622812431914202
$ws.Range("K5").Value = 0.2977848790824405
$ws.Range("L5").Value = -0.05749781206784608
$ws.Range("M5").Value = 0.07471037506210082
$ws.Range("N5").Value = 0.2627015921725792
$ws.Range("O5").Value = 0.2515737508594499
$ws.Range("B7").Value = -0.04502025361393325
$ws.Range("C7").Value = -0.01940981040636119
$ws.Range("D7").Value = 0.2858119586419813
$ws.Range("E7").Value = 0.9628637223458245
$ws.Range("H7").Value = 0.1850100977435571
$ws.Range("I7").Value = 0.5389219049867383
$ws.Range("J7").Value = 0.6093980989173357
$ws.Range("K7").Value = 0.4163707152817249
$ws.Range("L7").Value = -0.01143832450566349
$ws.Range("M7").Value = 0.2112297758292304
$ws.Range("N7").Value = 0.2815185517408797
$ws.Range("O7").Value = 0.3159095033480881
$ws.Range("B8").Value = -0.5352988047422513
$ws.Range("C8").Value = 0.1686289861464561
$ws.Range("D8").Value = 0.3208449815547014
$ws.Range("E8").Value = 0.1878429442749909
$ws.Range("G8").Value = 0.1850100977435571
$ws.Range("I8").Value = -0.2213412639926016
$ws.Range("J8").Value = 0.2664628392334352
$ws.Range("K8").Value = -0.1960758592444096
$ws.Range("L8").Value = -0.347064490796336
$ws.Range("M8").Value = -0.2171218224399398
$ws.Range("N8").Value = 0.9352104838129371
$ws.Range("O8").Value = -0.3576867119968038
$ws.Range("B9").Value = 0.1490515468083845
$ws.Range("C9").Value = 0.1051672521303078
$ws.Range("D9").Value = -0.03615914419776285
$ws.Range("E9").Value = 0.4840581834101727
$ws.Range("G9").Value = 0.5389219049867383
$ws.Range("H9").Value = -0.2213412639926016
$ws.Range("J9").Value = 0.2862144791013954
$ws.Range("K9").Value = 0.7884343634426342
$ws.Range("L9").Value = 0.08087592040379925
$ws.Range("M9").Value = 0.6644814596677858
$ws.Range("N9").Value = -0.1612105954626258
$ws.Range("O9").Value = 0.8915350258456595
$ws.Range("B10").Value = 0.1371387993843627
$ws.Range("C10").Value = -0.01637347079980878
$ws.Range("D10").Value = 0.8042684740305435
$ws.Range("E10").Value = 0.622812431914202
$ws.Range("G10").Value = 0.6093980989173357
$ws.Range("H10").Value = 0.2664628392334352
$ws.Range("I10").Value = 0.2862144791013954
$ws.Range("K10").Value = 0.3730356032872126
$ws.Range("L10").Value = -0.005120395540669093
$ws.Range("M10").Value = 0.2336922242139398
$ws.Range("N10").Value = 0.3370899062620499
$ws.Range("O10").Value = 0.1816867800423102
$ws.Range("B11").Value = 0.200626606284561
$ws.Range("C11").Value = -0.07889299984295564
$ws.Range("D11").Value = 0.1010925659905756
$ws.Range("E11").Value = 0.2977848790824405
$ws.Range("G11").Value = 0.4163707152817249
$ws.Range("H11").Value = -0.1960758592444096
$ws.Range("I11").Value = 0.7884343634426342
$ws.Range("J11").Value = 0.3730356032872126
$ws.Range("L11").Value = 0.2273945170531899
$ws.Range("M11").Value = 0.9569216284928839
$ws.Range("N11").Value = -0.07720906754130509
$ws.Range("O11").Value = 0.8213573610483902
$ws.Range("B12").Value = 0.1079330116962531
$ws.Range("C12").Value = 0.1237571271336534
$ws.Range("D12").Value = 0.1766723565517075
$ws.Range("E12").Value = -0.05749781206784608
$ws.Range("G12").Value = -0.01143832450566349
$ws.Range("H12").Value = -0.347064490796336
$ws.Range("I12").Value = 0.08087592040379925
$ws.Range("J12").Value = -0.005120395540669093
$ws.Range("K12").Value = 0.2273945170531899
$ws.Range("M12").Value = 0.266688302155974
$ws.Range("N12").Value = -0.3127709963931993
$ws.Range("O12").Value = 0.044570176597717
$ws.Range("B13").Value = 0.129293823016267
$ws.Range("C13").Value = -0.03019904530512181
$ws.Range("D13").Value = 0.06202386746339025
$ws.Range("E13").Value = 0.07471037506210082
$ws.Range("G13").Value = 0.2112297758292304
$ws.Range("H13").Value = -0.2171218224399398
$ws.Range("I13").Value = 0.6644814596677858
$ws.Range("J13").Value = 0.2336922242139398
$ws.Range("K13").Value = 0.9569216284928839
$ws.Range("L13").Value = 0.266688302155974
$ws.Range("N13").Value = -0.1188419959597674
$ws.Range("O13").Value = 0.7602727279385658
$ws.Range("B14").Value = -0.5502255047420385
$ws.Range("C14").Value = 0.1411896716252118
$ws.Range("D14").Value = 0.2608281976605054
$ws.Range("E14").Value = 0.2627015921725792
$ws.Range("G14").Value = 0.2815185517408797
$ws.Range("H14").Value = 0.9352104838129371
$ws.Range("I14").Value = -0.1612105954626258
$ws.Range("J14").Value = 0.3370899062620499
$ws.Range("K14").Value = -0.07720906754130509
$ws.Range("L14").Value = -0.3127709963931993
$ws.Range("M14").Value = -0.1188419959597674
$ws.Range("O14").Value = -0.3259720137548081
$ws.Range("B15").Value = 0.2537685189959536
$ws.Range("C15").Value = 0.02761189664673756
$ws.Range("D15").Value = -0.09255504229564779
$ws.Range("E15").Value = 0.2515737508594499
$ws.Range("G15").Value = 0.3159095033480881
$ws.Range("H15").Value = -0.3576867119968038
$ws.Range("I15").Value = 0.8915350258456595
$ws.Range("J15").Value = 0.1816867800423102
$ws.Range("K15").Value = 0.8213573610483902
$ws.Range("L15").Value = 0.044570176597717
$ws.Range("M15").Value = 0.7602727279385658
$ws.Range("N15").Value = -0.3259720137548081
